$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text for the new column G (set first so shared-string table order
# matches the author's: "Trajectory Length / Examinations" before
# "Trajectory Length (aka movements)")
$ws.Range("G1").Value = "Trajectory Length / Examinations"
# Match the wrap-text header style used by the other row-1 header cells
$ws.Range("G1").WrapText = $true

# Column A's header text was reworded
$ws.Range("A1").Value = "Trajectory Length (aka movements)"

# New column G: Trajectory Length / Examinations = A / C, for every data row.
# Enter the first row's formula standalone, then fill the remaining rows so the
# growth pattern mirrors how the author built it (single cell, then a fill
# down of the rest).
$ws.Range("G2").Formula = "=A2/C2"
$ws.Range("G3:G51").Formula = "=A3/C3"

# Give the new column a sensible custom width (matches column B's width)
$ws.Columns.Item(7).ColumnWidth = 16.83

# View state tweaks to mirror the author's final window/selection position
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 2
$ws.Range("I36").Select()
